$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every value in column C (rows 2-252, the "Fitness" column) with 7534,
# per the correction made to the SA algorithm / 746 logs.
$ws.Range("C2:C252").Value = 7534
